# 09-BulkUploadMachineGroup.xlsx - refresh the Machine Group / Machine Type
# lookup table with the current raw-material / machine-group reference data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new reference table only needs 16 data rows (was 20) - drop the
# now-unused trailing rows 17-20 entirely (incl. their row formatting).
$ws.Range("A17:B20").EntireRow.Delete() | Out-Null

# Full replacement data set: col A = Machine Group, col B = Machine Type.
$data = @(
    @("Machine Group", "Machine Type"),
    @("TCM_GROUP_1", "Tube cutting manual"),
    @("TCA_GROUP_1", "Pedrazzoli tube cutting"),
    @("PCM_GROUP_1", "Plasma"),
    @("TPP_GROUP_1", "Trump"),
    @("TBH_GROUP_1", "Horizontal tube bending"),
    @("TBV_GROUP_1", "Vertical tube bending"),
    @("SMC_GROUP_1", "Shearing"),
    @("PBM_GROUP_1", "Sheet bending"),
    @("NMC_GROUP_1", "Notching machine"),
    @("MPM_GROUP_1", "Mechanical press machine"),
    @("HPM_GROUP_1", "hydraullic press machine"),
    @("ssm_group_1", "Speedy seamer"),
    @("DMC_GROUP_1", "Drilling machine"),
    @("hpm_group_1", "hydraullic press machine"),
    @("TCM_GROUP_2", "Vertical band saw")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value2 = $data[$i][0]
    $ws.Cells.Item($row, 2).Value2 = $data[$i][1]
}

# Column widths were refit for the new (shorter) labels in col A and the
# new (longer) labels in col B.
$ws.Columns.Item(1).ColumnWidth = 13.830729166666666
$ws.Columns.Item(2).ColumnWidth = 22.330729166666668

# Selection moved off the now-deleted row 3 area down to the new first
# blank row (17) ready for further data entry.
$ws.Range("A17:B17").Select() | Out-Null
